# Applies the scheduled-runner profit recalculation update to Jenova_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 57: Quit Your Jib-jab / Gold Needle
$ws.Range("H57").Value = 52072.5
$ws.Range("J57").Value = 52072.5
$ws.Range("L57").Value = 156217.5
$ws.Range("N57").Value = -157215.5

# Row 61: Not Taking No for an Answer / Mega-Potion of Strength
$ws.Range("H61").Value = 31.4
$ws.Range("I61").Value = 31.4
$ws.Range("K61").Value = 94.19999999999999
$ws.Range("M61").Value = 77.80000000000001

# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 13894121
$ws.Range("I62").Value = 62501250
$ws.Range("K62").Value = 62501250
$ws.Range("M62").Value = -62500626

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 13894121
$ws.Range("I65").Value = 62501250
$ws.Range("K65").Value = 312506250
$ws.Range("M65").Value = -312503130

# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 46227.434
$ws.Range("I107").Value = 55592
$ws.Range("J107").Value = 1745.75
$ws.Range("K107").Value = 55592
$ws.Range("L107").Value = 1745.75
$ws.Range("M107").Value = -53672
$ws.Range("N107").Value = -5585.75

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 3831
$ws.Range("I116").Value = 3817.5
$ws.Range("J116").Value = 3851.25
$ws.Range("K116").Value = 3817.5
$ws.Range("L116").Value = 3851.25
$ws.Range("M116").Value = -375.5
$ws.Range("N116").Value = -10735.25

# Row 123: Nearly Bare / Gaja Grimoire
$ws.Range("H123").Value = 72850.57000000001
$ws.Range("J123").Value = 72850.57000000001
$ws.Range("L123").Value = 72850.57000000001
$ws.Range("N123").Value = -82650.57000000001

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2538.0833
$ws.Range("I137").Value = 1175.9584
$ws.Range("J137").Value = 5262.3335
$ws.Range("K137").Value = 3527.8752
$ws.Range("L137").Value = 15787.0005
$ws.Range("M137").Value = -977.8751999999999
$ws.Range("N137").Value = -20887.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 2811.4443
$ws.Range("I32").Value = 2811.4443
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2811.4443
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2524.4443
$ws.Range("N32").ClearContents()

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 6520.5557
$ws.Range("I45").Value = 1787.1
$ws.Range("K45").Value = 1787.1
$ws.Range("M45").Value = -1410.1

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 4498.609
$ws.Range("I61").Value = 2872.6875
$ws.Range("K61").Value = 2872.6875
$ws.Range("M61").Value = -2660.6875

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 13381.56
$ws.Range("I74").Value = 15027.95
$ws.Range("K74").Value = 15027.95
$ws.Range("M74").Value = -14153.95

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 13381.56
$ws.Range("I77").Value = 15027.95
$ws.Range("K77").Value = 75139.75
$ws.Range("M77").Value = -70771.75

# Row 80: A Squire to Inspire / Titanium Hoplon
$ws.Range("H80").Value = 64000
$ws.Range("J80").Value = 64000
$ws.Range("L80").Value = 64000
$ws.Range("N80").Value = -65996

# Row 83: All's Fair in Highborn Assassination (L) / Titanium Hoplon
$ws.Range("H83").Value = 64000
$ws.Range("J83").Value = 64000
$ws.Range("L83").Value = 192000
$ws.Range("N83").Value = -201984

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 153806.12
$ws.Range("I110").Value = 168792.97
$ws.Range("K110").Value = 168792.97
$ws.Range("M110").Value = -166747.97

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 9854.695
$ws.Range("I132").Value = 3606
$ws.Range("J132").Value = 14661.385
$ws.Range("K132").Value = 10818
$ws.Range("L132").Value = 43984.155
$ws.Range("M132").Value = -8288
$ws.Range("N132").Value = -49044.155

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4498.609
$ws.Range("I136").Value = 2872.6875
$ws.Range("K136").Value = 8618.0625
$ws.Range("M136").Value = -6068.0625

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2113411.8
$ws.Range("I86").Value = 1386401
$ws.Range("J86").Value = 4003639.5
$ws.Range("K86").Value = 1386401
$ws.Range("L86").Value = 4003639.5
$ws.Range("M86").Value = -1385278
$ws.Range("N86").Value = -4005885.5

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2113411.8
$ws.Range("I89").Value = 1386401
$ws.Range("J89").Value = 4003639.5
$ws.Range("K89").Value = 6932005
$ws.Range("L89").Value = 20018197.5
$ws.Range("M89").Value = -6926389
$ws.Range("N89").Value = -20029429.5

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 76159.42999999999
$ws.Range("I134").Value = 3359.111
$ws.Range("J134").Value = 207200
$ws.Range("K134").Value = 10077.333
$ws.Range("L134").Value = 621600
$ws.Range("M134").Value = -7542.332999999999
$ws.Range("N134").Value = -626670

$ws = $wb.Worksheets.Item("CRP")
# Row 52: Spin It Like You Mean It / Mahogany Spinning Wheel
$ws.Range("H52").Value = 67783.75
$ws.Range("I52").Value = 65625
$ws.Range("J52").Value = 69942.5
$ws.Range("K52").Value = 65625
$ws.Range("L52").Value = 69942.5
$ws.Range("M52").Value = -65331
$ws.Range("N52").Value = -70530.5

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 297462.72
$ws.Range("I58").Value = 528381.2
$ws.Range("J58").Value = 4965.933
$ws.Range("K58").Value = 528381.2
$ws.Range("L58").Value = 4965.933
$ws.Range("M58").Value = -528178.2
$ws.Range("N58").Value = -5371.933

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 3968.2144
$ws.Range("I99").Value = 2955.35
$ws.Range("K99").Value = 2955.35
$ws.Range("M99").Value = -1457.35

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 1354.9656
$ws.Range("I107").Value = 675.2143
$ws.Range("J107").Value = 1989.4
$ws.Range("K107").Value = 675.2143
$ws.Range("L107").Value = 1989.4
$ws.Range("M107").Value = 1244.7857
$ws.Range("N107").Value = -5829.4

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 3968.2144
$ws.Range("I126").Value = 2955.35
$ws.Range("K126").Value = 8866.049999999999
$ws.Range("M126").Value = -6396.049999999999

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 4249.909
$ws.Range("I132").Value = 4074.9
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 12224.7
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -9694.700000000001
$ws.Range("N132").Value = -23060

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 403714.6
$ws.Range("I134").Value = 3336.2354
$ws.Range("J134").Value = 1254518.6
$ws.Range("K134").Value = 10008.7062
$ws.Range("L134").Value = 3763555.8
$ws.Range("M134").Value = -7473.706200000001
$ws.Range("N134").Value = -3768625.8

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 297462.72
$ws.Range("I136").Value = 528381.2
$ws.Range("J136").Value = 4965.933
$ws.Range("K136").Value = 1585143.6
$ws.Range("L136").Value = 14897.799
$ws.Range("M136").Value = -1582593.6
$ws.Range("N136").Value = -19997.799

# Row 139: Weaving a Path / Acacia Spinning Wheel
$ws.Range("H139").Value = 93293.336
$ws.Range("J139").Value = 93293.336
$ws.Range("L139").Value = 93293.336
$ws.Range("N139").Value = -103573.336

$ws = $wb.Worksheets.Item("CUL")
# Row 38: Pretty as a Picture / Dark Vinegar
$ws.Range("H38").Value = 36.5
$ws.Range("I38").Value = 52.5
$ws.Range("J38").Value = 31.166666
$ws.Range("K38").Value = 157.5
$ws.Range("L38").Value = 93.49999800000001
$ws.Range("M38").Value = 189.5
$ws.Range("N38").Value = -787.499998

# Row 56: Culture Club / Crowned Pie
$ws.Range("H56").Value = 6852.4707
$ws.Range("I56").Value = 6852.4707
$ws.Range("K56").Value = 6852.4707
$ws.Range("M56").Value = -6322.4707

# Row 80: Saucy for a Suitor / Hollandaise Sauce
$ws.Range("H80").Value = 1067.3334
$ws.Range("I80").Value = 1101
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 3303
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -2367
$ws.Range("N80").Value = -4872

# Row 83: Saved by the Sauce (L) / Hollandaise Sauce
$ws.Range("H83").Value = 1067.3334
$ws.Range("I83").Value = 1101
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 9909
$ws.Range("L83").Value = 9000
$ws.Range("M83").Value = -5229
$ws.Range("N83").Value = -18360

$ws = $wb.Worksheets.Item("GSM")
# Row 51: When We Were Blings / Mythril Ear Cuffs
$ws.Range("H51").Value = 57000
$ws.Range("I51").Value = 50000
$ws.Range("J51").Value = 59333.332
$ws.Range("K51").Value = 50000
$ws.Range("L51").Value = 59333.332
$ws.Range("N51").Value = -60351.332
$ws.Range("M51").Value = -49491

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 3764.0908
$ws.Range("I102").Value = 3458.2307
$ws.Range("K102").Value = 3458.2307
$ws.Range("M102").Value = -1836.2307

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 508231
$ws.Range("I113").Value = 723907.6
$ws.Range("K113").Value = 723907.6
$ws.Range("M113").Value = -721737.6

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 264978.84
$ws.Range("J132").Value = 36528.3
$ws.Range("L132").Value = 109584.9
$ws.Range("N132").Value = -114644.9

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 746.36365
$ws.Range("I22").Value = 742.6316
$ws.Range("J22").Value = 770
$ws.Range("K22").Value = 742.6316
$ws.Range("L22").Value = 770
$ws.Range("M22").Value = -447.6316
$ws.Range("N22").Value = -1360

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 746.36365
$ws.Range("I27").Value = 742.6316
$ws.Range("J27").Value = 770
$ws.Range("K27").Value = 742.6316
$ws.Range("L27").Value = 770
$ws.Range("M27").Value = -635.6316
$ws.Range("N27").Value = -984

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 3742.8
$ws.Range("I46").Value = 2839
$ws.Range("K46").Value = 2839
$ws.Range("M46").Value = -2651

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 441.4138
$ws.Range("J55").Value = 1078
$ws.Range("L55").Value = 1078
$ws.Range("N55").Value = -1424

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 9288.223
$ws.Range("I61").Value = 7902
$ws.Range("J61").Value = 9684.286
$ws.Range("K61").Value = 7902
$ws.Range("L61").Value = 9684.286
$ws.Range("M61").Value = -7700
$ws.Range("N61").Value = -10088.286

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 3502.1667
$ws.Range("I93").Value = 3862.2
$ws.Range("K93").Value = 3862.2
$ws.Range("M93").Value = -2614.2

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 9288.223
$ws.Range("I113").Value = 7902
$ws.Range("J113").Value = 9684.286
$ws.Range("K113").Value = 7902
$ws.Range("L113").Value = 9684.286
$ws.Range("M113").Value = -5732
$ws.Range("N113").Value = -14024.286

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 6724.7334
$ws.Range("J132").Value = 9714
$ws.Range("L132").Value = 29142
$ws.Range("N132").Value = -34202

# Row 139: Giving Gatherers Their Gear / Gomphotherium Doublet of Gathering
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("WVR")
# Row 98: Apparent Apparel / Kudzu Tunic of Striking
$ws.Range("H98").Value = 82580
$ws.Range("J98").Value = 82580
$ws.Range("L98").Value = 82580
$ws.Range("N98").Value = -88570

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 1186.28
$ws.Range("I113").Value = 1098.5883
$ws.Range("K113").Value = 3295.7649
$ws.Range("M113").Value = -1125.7649
